$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.062.67"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "2.303.65"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'302.14"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "'99.30"
$ws.Range("E6").Value = "  +5.63%  "
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.507"
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("D10").Value = "'34.37"
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").Value = "'49.09"
$ws.Range("E12").Value = "  +3.20%  "
$ws.Range("E13").Value = "  +4.25%  "
$ws.Range("D14").Value = "'17.84"
$ws.Range("E14").Value = "  +17.00%  "
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "2.664.04"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "2.287.52"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("E18").Value = "  +4.35%  "
$ws.Range("D19").Value = "42.951.27"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("E20").Value = "  +8.64%  "
$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("D22").Value = "'6.10"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").Value = "'67.79"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("D24").Value = "'237.09"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("E25").Value = "  +13.21%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.46"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "'24.51"
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("D29").Value = "'168.11"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").Value = "'2.08"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("D31").Value = "'33.85"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").Value = "'9.17"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("E35").Value = "  +4.51%  "
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").Value = "'17.01"
$ws.Range("E37").Value = "  +6.83%  "
$ws.Range("D38").Value = "'0.0699"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("E40").Value = "  +4.78%  "
$ws.Range("D41").Value = "'2.81"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").Value = "2.002.50"
$ws.Range("E44").Value = "  +2.88%  "
$ws.Range("E45").Value = "  +2.81%  "
$ws.Range("D46").Value = "'10.09"
$ws.Range("E46").Value = "  +5.72%  "
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("D48").Value = "'2.86"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("D49").Value = "'55.35"
$ws.Range("E49").Value = "  +6.16%  "
$ws.Range("D50").Value = "2.529.29"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("E51").Value = "  +2.46%  "
